# Add fishing mini game
#  - new eventId 104 ("Fishing mini game") in the MapEvent lookup sheet
#  - five tiles on MapConfig switched from the generic "EARN" (100) event
#    to the new "FISH" (104) event
#  - view/selection state updated to match where the author was working

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MapConfig")
$ws2 = $wb.Worksheets.Item("MapEvent")

# --- MapEvent: describe the new event (eventId 104 already exists as a row,
# only its description text changes) -----------------------------------
# Set this FIRST so "Fishing mini game" becomes the earlier of the two new
# shared strings (matches the source order in the workbook).
$ws2.Range("C8").Value = "Fishing mini game"

# --- MapConfig: retarget the five tiles onto the fishing event ---------
$fishRows = @(6, 14, 23, 34, 42)
foreach ($r in $fishRows) {
    $ws1.Range("F$r").Value = 104
    $ws1.Range("I$r").Value = "FISH"
}

# --- view state -----------------------------------------------------------
# MapEvent's selection moves to C12, but MapConfig stays the active tab.
$ws2.Range("C12").Select()
$ws1.Activate()

$app = $wb.Application
$win = $app.ActiveWindow
try { $win.ScrollRow = 19 } catch {}
try { $win.ScrollColumn = 1 } catch {}

$ws1.Range("N38").Select()
